$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. New "LDCP"/"STCP" coprocessor instructions + their bit-diagram.
#    New shared strings must be created in this exact order so they
#    land on the same sharedStrings.xml indices as the target file:
#    LDCP, Cp-sel, CR=coprocessor register, CR-4, CR-3, CR-2, CR-1,
#    CR-0, STCP.
# ------------------------------------------------------------------

# LDCP row (row 44) - first use of "LDCP" creates the new string.
$ws.Range("B44").Value = "LDCP"

# Bit-diagram header for the new Cp-sel format (row 14).
$ws.Range("G14").Value = "Cp-sel"
$ws.Range("H14").Value = "CR=coprocessor register"

# Bit numbering row (row 15), mirrors the other "bits" rows (2/7/11).
$ws.Range("G15").Value = "bits"
$ws.Range("H15").Value = 7
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = 4
$ws.Range("L15").Value = 3
$ws.Range("M15").Value = 2
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0

# Field values row (row 16): coprocessor register select + rD x3.
$ws.Range("H16").Value = "CR-4"
$ws.Range("I16").Value = "CR-3"
$ws.Range("J16").Value = "CR-2"
$ws.Range("K16").Value = "CR-1"
$ws.Range("L16").Value = "CR-0"
$ws.Range("M16").Value = "rD"
$ws.Range("N16").Value = "rD"
$ws.Range("O16").Value = "rD"

# Finish the LDCP row, and add the STCP row right below it.
$ws.Range("C44").Value = "Cp-sel"
$ws.Range("B45").Value = "STCP"
$ws.Range("C45").Value = "Cp-sel"

# ------------------------------------------------------------------
# 2. The opcode quick-reference list (column G/H, previously rows
#    17-42) shifts down two rows to 19-44 to make room for the new
#    bit diagram above, and gains one more entry (0x23/0x24 pair).
#    Rows themselves are not inserted (column A-D data keeps its
#    row numbers) - only this side list's cell contents move.
# ------------------------------------------------------------------
$ws.Range("G17").ClearContents()
$ws.Range("G18").ClearContents()
$ws.Range("H18").ClearContents()

$ws.Range("G19").Value = "Spec-Opc"
$ws.Range("G20").Value = "0x00"
$ws.Range("H20").Value = "RET"
$ws.Range("G21").Value = "0x01"
$ws.Range("G22").Value = "0x02"
$ws.Range("G23").Value = "0x03"
$ws.Range("G24").Value = "0x04"
$ws.Range("G25").Value = "0x05"
$ws.Range("G26").Value = "0x06"
$ws.Range("G27").Value = "0x07"
$ws.Range("G28").Value = "0x08"
$ws.Range("G29").Value = "0x09"
$ws.Range("G30").Value = "0x10"
$ws.Range("G31").Value = "0x11"
$ws.Range("G32").Value = "0x12"
$ws.Range("G33").Value = "0x13"
$ws.Range("G34").Value = "0x14"
$ws.Range("G35").Value = "0x15"
$ws.Range("G36").Value = "0x16"
$ws.Range("G37").Value = "0x17"
$ws.Range("G38").Value = "0x18"
$ws.Range("G39").Value = "0x19"
$ws.Range("G40").Value = "0x20"
$ws.Range("G41").Value = "0x21"
$ws.Range("G42").Value = "0x22"
$ws.Range("G43").Value = "0x23"
$ws.Range("G44").Value = "0x24"

# ------------------------------------------------------------------
# 3. Swap the rS/rD field order in the two "X" (modifier) bit rows.
# ------------------------------------------------------------------
$ws.Range("J3").Value = "rS"
$ws.Range("K3").Value = "rS"
$ws.Range("L3").Value = "rS"
$ws.Range("M3").Value = "rD"
$ws.Range("N3").Value = "rD"
$ws.Range("O3").Value = "rD"

$ws.Range("J12").Value = "shift-2"
$ws.Range("K12").Value = "shift-1"
$ws.Range("L12").Value = "Shift-0"
$ws.Range("M12").Value = "rD"
$ws.Range("N12").Value = "rD"
$ws.Range("O12").Value = "rD"

# ------------------------------------------------------------------
# 4. Selection / scroll position: no longer scrolled to A30, the
#    active cell moves from I41 to B46.
# ------------------------------------------------------------------
$ws.Range("B46").Select()
